# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Femacal de La Calera" / Sandia at row 514,
# shifting the existing rows 514-525 down to 515-526.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 514 (shift existing rows down).
$ws.Rows.Item(514).Insert(-4121)  # -4121 = xlShiftDown

# Populate the newly inserted row 514 with the new data point.
$ws.Range("A514").Value = 3
$ws.Range("B514").Value = "Femacal de La Calera"
$ws.Range("C514").Value = "Coquimbo"
$ws.Range("D514").Value2 = 44890
$ws.Range("E514").Value = 5
$ws.Range("F514").Value = 100112028
$ws.Range("G514").Value = "Sandia"
$ws.Range("H514").Value = "Sin especificar"
$ws.Range("I514").Value = "Primera"
$ws.Range("J514").Value = 290
$ws.Range("K514").Value = 800
$ws.Range("L514").Value = 850
$ws.Range("M514").Value = 831
$ws.Range("N514").Value = '$/kilo (volumen en unidades)'
$ws.Range("O514").Value = "Perú"
$ws.Range("P514").Value = 831
$ws.Range("Q514").Value = 1
$ws.Range("R514").Value = "Hortaliza"
